$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row-level Price (D) / Volume(1h) (E) refreshes ---
$ws.Range("D2").Value = "69.082.37"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "3.818.46"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'629.87"
$ws.Range("E5").Value = "  +5.44%  "
$ws.Range("D6").Value = "'165.33"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").Value = "3.815.75"
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("E10").Value = "  +2.65%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").Value = "'6.63"
$ws.Range("E12").Value = "  +3.56%  "
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "'36.02"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").Value = "4.459.11"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "3.969.17"
$ws.Range("E16").Value = "  +4.79%  "
$ws.Range("D17").Value = "69.053.99"
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").Value = "'18.03"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").Value = "'7.12"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'465.65"
$ws.Range("E21").Value = "  +1.18%  "
$ws.Range("D22").Value = "'9.65"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("E24").Value = "  +4.15%  "
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").Value = "'11.96"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +3.17%  "
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "3.967.38"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").Value = "'2.69"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("D32").Value = "'2.23"
$ws.Range("E32").Value = "  +2.10%  "
$ws.Range("D33").Value = "'7.29"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").Value = "'29.20"
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("E37").Value = "  +3.24%  "
$ws.Range("E38").Value = "  +8.19%  "
$ws.Range("D39").Value = "'3.42"
$ws.Range("E39").Value = "  +5.77%  "
$ws.Range("E40").Value = "  +3.38%  "
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D44").Value = "'157.59"
$ws.Range("E44").Value = "  +4.11%  "
$ws.Range("E45").Value = "  +5.54%  "
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").Value = "'46.89"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D50").Value = "'42.47"
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("E51").Value = "  +13.29%  "

# --- Rows 48/49 swapped: Cosmos and Stacks traded ranking positions ---
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.90"
$ws.Range("E48").Value = "  +3.25%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'8.45"
$ws.Range("E49").Value = "  +1.75%  "
